# Updated cryptos list (prices + 1h volume %) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.541.45'
$ws.Range('E2').Value = '  +0.89%  '

$ws.Range('D3').Value = '1.880.86'
$ws.Range('E3').Value = '  +1.20%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7178'
$ws.Range('E5').Value = '  +2.17%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.00'
$ws.Range('E6').Value = '  +1.67%  '

$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07953'
$ws.Range('E8').Value = '  +0.38%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3100'
$ws.Range('E9').Value = '  +2.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.35'
$ws.Range('E10').Value = '  +3.61%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08276'
$ws.Range('E11').Value = '  +1.36%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.887.90'
$ws.Range('E12').Value = '  +0.81%  '

$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7295'
$ws.Range('E13').Value = '  +3.17%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.286'
$ws.Range('E14').Value = '  +1.37%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.23'
$ws.Range('E15').Value = '  +1.92%  '

$ws.Range('D16').Value = '29.538.92'
$ws.Range('E16').Value = '  +0.65%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.909'
$ws.Range('E17').Value = '  +1.66%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '246.33'
$ws.Range('E18').Value = '  +3.74%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007880'
$ws.Range('E19').Value = '  +0.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.35'
$ws.Range('E20').Value = '  +0.99%  '

$ws.Range('D21').Value = '2.120.60'
$ws.Range('E21').Value = '  -0.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.069'
$ws.Range('E22').Value = '  +6.63%  '

$ws.Range('E23').Value = '  +0.21%  '

$ws.Range('E24').Value = '  +0.19%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1606'
$ws.Range('E25').Value = '  +13.10%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.74'
$ws.Range('E26').Value = '  +0.76%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.049'
$ws.Range('E27').Value = '  +1.63%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.34'
$ws.Range('E28').Value = '  +1.49%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.357'
$ws.Range('E29').Value = '  -3.09%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.495'
$ws.Range('E30').Value = '  +1.14%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.403'
$ws.Range('E31').Value = '  +2.49%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.109'
$ws.Range('E32').Value = '  +1.90%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05212'
$ws.Range('E33').Value = '  +0.69%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.946'
$ws.Range('E34').Value = '  +1.85%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.200'
$ws.Range('E35').Value = '  +1.76%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7256'
$ws.Range('E36').Value = '  +2.34%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.685'
$ws.Range('E37').Value = '  +0.24%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01874'
$ws.Range('E38').Value = '  +1.30%  '

$ws.Range('D39').Value = '1.203.45'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.708'
$ws.Range('E40').Value = '  +0.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9108'
$ws.Range('E41').Value = '  -1.09%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.166'
$ws.Range('E42').Value = '  +3.45%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '73.66'
$ws.Range('E43').Value = '  +4.63%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  +0.21%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.43'
$ws.Range('E45').Value = '  -0.45%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5301'
$ws.Range('E46').Value = '  -0.20%  '

$ws.Range('D47').Value = '2.016.18'
$ws.Range('E47').Value = '  -0.61%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.798'
$ws.Range('E48').Value = '  +2.98%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.936'
$ws.Range('E49').Value = '  +9.46%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.314'
$ws.Range('E50').Value = '  +1.40%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4329'
$ws.Range('E51').Value = '  +1.97%  '
